$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.375.77"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.36"
$ws.Range("E3").Value = "  +5.39%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.99"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5671"
$ws.Range("E7").Value = "  +17.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3891"
$ws.Range("E8").Value = "  +12.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07624"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.11"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  +8.83%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +7.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.258"
$ws.Range("E14").Value = "  +6.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.811.76"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("E16").Value = "  +7.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.05"
$ws.Range("E17").Value = "  +6.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001077"
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06485"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  +4.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.012"
$ws.Range("E22").Value = "  +5.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.390.11"
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.120"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.87"
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.77"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.435"
$ws.Range("E28").Value = "  +17.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.023.26"
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.77"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.176"
$ws.Range("E31").Value = "  +13.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("E32").Value = "  +13.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.783"
$ws.Range("E33").Value = "  +7.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.636"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.918"
$ws.Range("E35").Value = "  +20.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02320"
$ws.Range("E36").Value = "  +6.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2160"
$ws.Range("E37").Value = "  +8.73%  "
$ws.Range("E38").Value = "  +7.25%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06098"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6417"
$ws.Range("E40").Value = "  +7.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.041"
$ws.Range("E41").Value = "  +6.69%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.154"
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.380"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5991"
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.700"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.944"
$ws.Range("E49").Value = "  +5.90%  "
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06849"
$ws.Range("E51").Value = "  +3.36%  "
